$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: correct the student's name / gender / parent email ---
$ws.Range("A2").Value = "Aremu"
$ws.Range("B2").Value = "Asade"
$ws.Range("C2").Value = "Openiyi"
$ws.Range("D2").Value = "M"
$ws.Range("Q2").Value = "AremuAsade@yahoo.com"

# --- Row 3: brand-new student record ---
$ws.Range("A3").Value = "Alade"
$ws.Range("B3").Value = "Abiola"
$ws.Range("C3").Value = "Babatunde"
$ws.Range("D3").Value = "F"
$ws.Range("E3").Value = "Staff"
$ws.Range("F3").Value = 38691
$ws.Range("G3").Value = 2021
$ws.Range("H3").Value = "lagos"
$ws.Range("I3").Value = "ajeromi ifelodun"
$ws.Range("J3").Value = "Christianity"
$ws.Range("K3").Value = "18 Jones Waribi"
$ws.Range("L3").Value = "Lagos"
$ws.Range("M3").Value = "Lagos"
$ws.Range("N3").Value = "Mr. ABDULLAHI"
$ws.Range("O3").Value = "ABDULLAHI"
$ws.Range("P3").Value = "Male"
$ws.Range("Q3").Value = "Olawale099@yahoo.com"
$ws.Range("R3").Value = 9089897766
$ws.Range("S3").Value = "Nigerian"
$ws.Range("T3").Value = "Lagos"
$ws.Range("U3").Value = "Lagos"
$ws.Range("V3").Value = "18 Jones Waribi"
$ws.Range("W3").Value = "Business"
$ws.Range("X3").Value = "lagos"
$ws.Range("Y3").Value = "ajeromi ifelodun"
$ws.Range("Z3").Value = "christian"

# --- Row 4: another brand-new student record ---
$ws.Range("A4").Value = "Olufunke"
$ws.Range("B4").Value = "Akinkunmi"
$ws.Range("C4").Value = "Olawale"
$ws.Range("D4").Value = "M"
$ws.Range("E4").Value = "Staff"
$ws.Range("F4").Value = 38691
$ws.Range("G4").Value = 2021
$ws.Range("H4").Value = "lagos"
$ws.Range("I4").Value = "ajeromi ifelodun"
$ws.Range("J4").Value = "Christianity"
$ws.Range("K4").Value = "18 Jones Waribi"
$ws.Range("L4").Value = "Lagos"
$ws.Range("M4").Value = "Lagos"
$ws.Range("N4").Value = "Mr. ABDULLAHI"
$ws.Range("O4").Value = "ABDULLAHI"
$ws.Range("P4").Value = "Male"
$ws.Range("Q4").Value = "tundetunapa@@yahoo.com"
$ws.Range("R4").Value = 9089897766
$ws.Range("S4").Value = "Nigerian"
$ws.Range("T4").Value = "Lagos"
$ws.Range("U4").Value = "Lagos"
$ws.Range("V4").Value = "18 Jones Waribi"
$ws.Range("W4").Value = "Business"
$ws.Range("X4").Value = "lagos"
$ws.Range("Y4").Value = "ajeromi ifelodun"
$ws.Range("Z4").Value = "christian"

# --- Hyperlinks for the parent e-mail addresses entered in Q2 / Q3 ---
$ws.Hyperlinks.Add($ws.Range("Q2"), "mailto:AremuAsade@yahoo.com")
$ws.Range("Q2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("Q3"), "mailto:Olawale099@yahoo.com")
$ws.Range("Q3").Style = "Hyperlink"

# --- View tweaks recorded in the workbook ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H13").Select()
